$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The J1-J4 4-pin connector (row 2) is swapped for a 2x4 nanofit connector that
# only uses J1,J2 - the old J3,J4 references move over to the Conn_01x03 connector.
$ws.Range("A2").Value = "J1,J2"
$ws.Range("B2").Value = "Conn_02x04"
$ws.Range("C2").Value = "1053141108"
$ws.Range("D2").Value = "Scrutineering:nanofit_02x04"
$ws.Range("E2").Value = 2

# Conn_01x03 (row 4) reference renumbered from J11,J12 to J3,J4
$ws.Range("A4").Value = "J3,J4"

# Conn_01x02 (row 5) reference renumbered from J13,J14 to J11,J11
$ws.Range("A5").Value = "J11,J11"

# Re-sort the BOM rows by the Value column, ascending, keeping the header row fixed
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B1:B11"))
$sortObj.SetRange($ws.Range("A1:E11"))
$sortObj.Header = 1
$sortObj.Apply()

$ws.Range("N11").Select()
